# PIXm webinar 2021 - update version references from Rev 2.1 to Rev 3.0 (slide 6)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# --- Title shape ("Significant changes from PIXm, Rev 2.1") ---------------
# Collapse the three runs into the text "Significant changes from PIXm, Rev 3.0"
# by trimming everything after "Significant changes from " and retyping the
# remainder; this keeps the first run's character formatting (incl. dirty="0")
# and produces a single merged run, matching the target markup.
$titleShape = $s.Shapes.Item(2)
$titleRange = $titleShape.TextFrame.TextRange
$titleKeepLen = 25 # Len("Significant changes from ")
$titleTail = $titleRange.Characters($titleKeepLen + 1, $titleRange.Length - $titleKeepLen)
$titleTail.Text = ""
$titleRange.Text = $titleRange.Text + "PIXm, Rev 3.0"

# --- Subtitle shape ("FHIR Implementation Guide instead of PDF - Rev. 2.1") -
# Drop the trailing " - Rev. 2.1" and split "of PDF" into its own run (this
# mirrors the target where "FHIR Implementation Guide instead " stays one run
# and "of PDF" becomes a second run with matching formatting).
$subShape = $s.Shapes.Item(3)
$subRange = $subShape.TextFrame.TextRange
$suffix = $subRange.Characters(41, 11) # " - Rev. 2.1"
$suffix.Text = ""
$ofPdf = $subRange.Characters(35, 6) # "of PDF"
$ofPdf.Font.Name = "Arial"

Write-Host "Updated slide 6 title and subtitle to Rev 3.0"
